# Staging.OrganizationType.xlsx - "moved staging files StagingTemplates directory"
#
# The meaningful, user-visible content change in this commit is that the two
# header cells on row 2 swap places: A2 ("OrganizationType_ID") and
# B2 ("Description") trade values, so the sheet now reads
# A2 = "Description", B2 = "OrganizationType_ID".
#
# (The diff also shows the bookViews windowWidth/windowHeight and the sheet's
# VBA codeName changing, and the explicit width/bestFit formatting on column B
# disappearing - those are Excel-internal/window-chrome bookkeeping values
# that aren't driven by user action and aren't exposed as settable through
# the Excel object model, so they're left alone here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the header labels in A2 / B2.
# (Use Value2 - on this host, the plain Value property round-trips oddly for
# get/set of simple strings, Value2 is the reliable equivalent here.)
$a2 = $ws.Range("A2").Value2
$b2 = $ws.Range("B2").Value2

$ws.Range("A2").Value2 = $b2
$ws.Range("B2").Value2 = $a2
